# Updated symbol list on Mon Dec 12 17:45:44 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'276.20"
$ws.Range("D3").Value = "'21.06"
$ws.Range("D4").Value = "'6.213"
$ws.Range("D5").Value = "'0.06188"
$ws.Range("D6").Value = "'3.579"
$ws.Range("D7").Value = "'1.518"
$ws.Range("D8").Value = "'6.551"
$ws.Range("D9").Value = "'0.8227"
$ws.Range("D10").Value = "'0.1643"
$ws.Range("D11").Value = "'0.08233"
$ws.Range("D12").Value = "'0.03441"
$ws.Range("D13").Value = "'0.03128"
$ws.Range("D14").Value = "'0.09130"
$ws.Range("D15").Value = "'3.769"
$ws.Range("D16").Value = "'0.001624"
$ws.Range("D17").Value = "'0.04701"
$ws.Range("D18").Value = "'0.006315"
$ws.Range("D19").Value = "'0.006141"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("D21").Value = "'0.0001500"
$ws.Range("D24").Value = "'0.01387"
$ws.Range("D25").Value = "'0.3282"
$ws.Range("D28").Value = "'0.0002737"
$ws.Range("D40").Value = "'0.04670"
$ws.Range("D41").Value = "'0.007061"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.004601"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1103"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").Value = "'0.01115"
$ws.Range("D45").Value = "'0.00006251"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D47").Value = "'0.8453"
$ws.Range("D49").Value = "'0.00001900"
$ws.Range("D50").Value = "'0.01240"
